# Applies the Summary / Experience / Education / Skills / Projects rewrite.
$d = $word.ActiveDocument

# Word's "manual line break" character, as produced by Shift+Enter.
# When assigned into a Range.Text string it serializes to an OOXML <w:br/>.
$nl = [char]11

# --- Summary ---------------------------------------------------------------
$d.Content.Find.Execute(
    "A student with a strong enthusiasm for technology and app development, leveraging Computer Science Engineering knowledge to drive innovation and growth in the tech industry through the application of academic foundations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dedicated and people-oriented professional with a background in education and experience in supporting and communicating with individuals, seeking a role in Human/Social Services or a related field. Strong interpersonal skills and a passion for helping others are key strengths, poised for growth in a dynamic environment focused on social assistance and administration.",
    2) | Out-Null

# --- Experience --------------------------------------------------------------
$d.Content.Find.Execute(
    "Developed an AI-powered resume web application using AI models, training and fine-tuning them to generate high-quality content. Optimized model output through prompt refinement, leveraging AI capabilities to drive innovative resume development. Demonstrated expertise in AI model training and deployment, yielding enhanced resume generation with measurable impact through successful model output optimization.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed an AI web-app for resume enhancement using LLM models, training them to identify required resume features and fine-tuning with targeted prompts for content generation, resulting in improved content creation efficiency.",
    2) | Out-Null

# --- Education ---------------------------------------------------------------
# Needs a line break inside the paragraph, so find the old sentence, then
# overwrite the matched Range's Text (rather than using Find's plain
# Replacement, which cannot carry a line-break character).
$rngEdu = $d.Content
$rngEdu.Find.Execute(
    "Bachelor of Engineering in Computer Science and Engineering, AVIT, May 2026, GPA 7.1. Relevant coursework includes computer science and engineering fundamentals.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$eduLines = @(
    "Bachelor of Engineering in Computer Science, AVIT, May 2026, GPA 7.1;",
    "Intermediate, Narayana Jr College, June 2022, GPA 8.9"
)
$rngEdu.Text = $eduLines -join $nl

# --- Skills -------------------------------------------------------------------
$rngSkills = $d.Content
$rngSkills.Find.Execute("C#, Node, Python", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$skillsLines = @(
    "Good with people, can talk and listen well but sometimes nervous,",
    "Know some billing and petty cash handling but not expert,",
    "Experience running group sessions for mental health and social skills,",
    "Can plan and do social activities but not very creative,",
    "",
    "Able to multitask but sometimes get overwhelmed,",
    "",
    "Basic computer skills like Microsoft Word and Excel, "
)
$rngSkills.Text = $skillsLines -join $nl

# --- Projects -----------------------------------------------------------------
$oldProjectLines = @(
    "Developed a QR scanner and generator web-app with TypeScript and Node, enabling secure sharing of encrypted information via PIN protection.",
    "Created Prediction PRO, a Solana-based crypto staking and price prediction app providing real-time forecasts.",
    "Designed a purchase order management system, allowing users to track orders from multiple apps, receive delivery reminders, and consolidate management in a single interface."
)
$rngProj = $d.Content
$rngProj.Find.Execute(
    ($oldProjectLines -join $nl),
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngProj.Text = "Developed and implemented a QR scanner and generator, leveraging TypeScript and Node.js for the backend. Spearheaded backend development, ensuring seamless functionality. Additionally, contributed to Prediction Pro, a full-stack application built with TypeScript, React, and PostgreSQL, demonstrating expertise in modern technologies and collaborative skills through frontend and database integration. These projects showcased versatility in tech stacks, full-stack development capabilities, and effective collaboration."
